$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (price + 1h volume refresh); rows 15/16 swap rank order.
# D-column prices are forced to text (leading apostrophe) so Excel's
# auto-type-detection does not reinterpret/reformat numeric-looking price
# strings (e.g. dropping trailing zeros or using thousands-dot grouping).
$ws.Range("D2").Value = "'24.611.80"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "'1.675.41"
$ws.Range("E3").Value = "  -1.79%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'313.35"
$ws.Range("E5").Value = "  -1.19%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "'0.3897"
$ws.Range("E7").Value = "  -3.21%  "
$ws.Range("D8").Value = "'0.3932"
$ws.Range("E8").Value = "  -3.32%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").Value = "'51.98"
$ws.Range("E10").Value = "  -3.27%  "
$ws.Range("D11").Value = "'1.389"
$ws.Range("E11").Value = "  -6.35%  "
$ws.Range("D12").Value = "'0.08619"
$ws.Range("E12").Value = "  -2.25%  "
$ws.Range("D13").Value = "'25.07"
$ws.Range("E13").Value = "  -4.73%  "
$ws.Range("D14").Value = "'7.290"
$ws.Range("E14").Value = "  -2.68%  "

# Rows 15-16: ShibaInu and Chainlink swap rank positions (ShibaInu moves to 15, Chainlink to 16)
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.00001317"
$ws.Range("E15").Value = "  -2.88%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'7.752"
$ws.Range("E16").Value = "  -4.29%  "

$ws.Range("D17").Value = "'1.690.44"
$ws.Range("E17").Value = "  -1.78%  "
$ws.Range("D18").Value = "'93.60"
$ws.Range("E18").Value = "  -3.29%  "
$ws.Range("D19").Value = "'0.07062"
$ws.Range("E19").Value = "  -1.38%  "
$ws.Range("D20").Value = "'20.55"
$ws.Range("E20").Value = "  -2.44%  "
$ws.Range("D21").Value = "'7.051"
$ws.Range("E21").Value = "  -2.75%  "
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").Value = "'13.99"
$ws.Range("E23").Value = "  -2.73%  "
$ws.Range("D24").Value = "'24.608.83"
$ws.Range("E24").Value = "  -1.17%  "
$ws.Range("D25").Value = "'2.368"
$ws.Range("E25").Value = "  +1.91%  "
$ws.Range("D26").Value = "'23.13"
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("D27").Value = "'2.725"
$ws.Range("E27").Value = "  -5.99%  "
$ws.Range("D28").Value = "'162.25"
$ws.Range("E28").Value = "  -2.53%  "
$ws.Range("D29").Value = "'5.819"
$ws.Range("E29").Value = "  -7.24%  "
$ws.Range("D30").Value = "'146.58"
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("D31").Value = "'7.984"
$ws.Range("E31").Value = "  -3.43%  "
$ws.Range("D32").Value = "'2.575"
$ws.Range("E32").Value = "  +15.32%  "
$ws.Range("D33").Value = "'1.874.76"
$ws.Range("E33").Value = "  -2.65%  "
$ws.Range("D34").Value = "'0.08362"
$ws.Range("E34").Value = "  -5.83%  "
$ws.Range("D35").Value = "'0.03039"
$ws.Range("E35").Value = "  -5.18%  "
$ws.Range("D36").Value = "'0.2817"
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("D37").Value = "'6.843"
$ws.Range("E37").Value = "  -4.92%  "
$ws.Range("D38").Value = "'0.9814"
$ws.Range("E38").Value = "  -4.33%  "
$ws.Range("D39").Value = "'0.09511"
$ws.Range("E39").Value = "  +2.24%  "
$ws.Range("D40").Value = "'1.542"
$ws.Range("E40").Value = "  +4.65%  "
$ws.Range("D41").Value = "'10.46"
$ws.Range("E41").Value = "  -3.62%  "
$ws.Range("D42").Value = "'0.7861"
$ws.Range("E42").Value = "  -7.30%  "
$ws.Range("D43").Value = "'13.51"
$ws.Range("E43").Value = "  -4.99%  "
$ws.Range("D44").Value = "'16.46"
$ws.Range("E44").Value = "  -6.10%  "
$ws.Range("D45").Value = "'0.7099"
$ws.Range("E45").Value = "  -4.82%  "
$ws.Range("D46").Value = "'2.553"
$ws.Range("E46").Value = "  -6.42%  "
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("D48").Value = "'0.08630"
$ws.Range("E48").Value = "  +3.09%  "
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").Value = "'1.319"
$ws.Range("E50").Value = "  -5.54%  "
$ws.Range("D51").Value = "'137.06"
$ws.Range("E51").Value = "  -3.69%  "
